$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain text so numeric-looking strings
# ("57.904.38", "0.386", "6.64", ...) are not auto-coerced to numbers by
# Excel's type inference on Range.Value assignment.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Cell value updates (crypto price / 1h-volume refresh) ---
$ws.Range("D2").Value = "57.904.38"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "3.121.72"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "525.69"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("D6").Value = "141.53"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.122.19"
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").Value = "0.386"
$ws.Range("E12").Value = "  +3.40%  "
$ws.Range("D13").Value = "3.658.98"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("E15").Value = "  +3.95%  "
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").Value = "57.997.87"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "3.116.92"
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "337.69"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("D25").Value = "66.84"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "0.0₃0934"
$ws.Range("E28").Value = "  +3.53%  "
$ws.Range("D29").Value = "6.64"
$ws.Range("E29").Value = "  +4.93%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "7.26"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("E32").Value = "  +3.18%  "
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "4.68"
$ws.Range("E35").Value = "  +5.12%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "154.13"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").Value = "6.14"
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("D38").Value = "27.11"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "3.161.70"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").Value = "0.688"
$ws.Range("E42").Value = "  +5.34%  "
$ws.Range("E43").Value = "  +10.69%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "37.08"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "2.304.88"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  +7.70%  "
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").Value = "6.04"
$ws.Range("E51").Value = "  +2.98%  "

# Restore column D to its original (default) style now that the text
# values are committed, so no stray number-format style lingers on the cells.
$ws.Range("D2:D51").Style = "Normal"

